$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the redundant header row ("cidade", "Casos confirmados", "Óbitos confirmados")
# at row 2 - rows below shift up by one.
$ws.Rows("2:2").Delete()

# Remove the trailing "outros estados" / "outros paises" rows, which after the
# shift above now sit at rows 37 and 38.
$ws.Rows("37:38").Delete()
